# Adapt the AHB-Diff header row to the new "<formatversion>" naming scheme,
# wrap the sheet's used range in an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row (A1:U1): "<label>_old" -> "<label>_FV2304", "<label>_new" -> "<label>_FV2310"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2304")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2310")
    }
}

# 2. Turn the used range A1:U55 into an Excel Table ("Table1") with an AutoFilter.
$rng = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $rng, $false, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3. Freeze the header row (split/freeze after row 1, topLeftCell A2).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
